$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.379.96"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "2.915.37"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.16"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.50"
$ws.Range("E6").Value = "  +0.79%  "

$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.08"
$ws.Range("E10").Value = "  -0.76%  "

$ws.Range("E11").Value = "  +3.60%  "

$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.89"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.83"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").Value = "3.373.34"
$ws.Range("E15").Value = "  +3.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.01"
$ws.Range("E16").Value = "  +6.07%  "

$ws.Range("D17").Value = "2.934.37"
$ws.Range("E17").Value = "  +4.69%  "

$ws.Range("D18").Value = "52.410.84"
$ws.Range("E18").Value = "  +1.23%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("E21").Value = "  +4.19%  "

$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.00"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.25"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("E26").Value = "  +7.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.83"
$ws.Range("E27").Value = "  +2.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.68"
$ws.Range("E29").Value = "  +2.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.68"
$ws.Range("E30").Value = "  +8.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.40"
$ws.Range("E31").Value = "  +13.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.01"
$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0981"
$ws.Range("E34").Value = "  +11.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.22"
$ws.Range("E35").Value = "  +1.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +6.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.93"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.08"
$ws.Range("E40").Value = "  +3.48%  "

$ws.Range("E41").Value = "  +14.13%  "

$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.47"
$ws.Range("E43").Value = "  +6.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.42"
$ws.Range("E44").Value = "  +1.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.60"
$ws.Range("E45").Value = "  +6.98%  "

$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("E47").Value = "  +4.50%  "

$ws.Range("D48").Value = "2.202.72"
$ws.Range("E48").Value = "  +4.00%  "

$ws.Range("E49").Value = "  +22.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0344"
$ws.Range("E50").Value = "  +11.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.970"
$ws.Range("E51").Value = "  +1.73%  "
